$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = 1
}

$ws.Range("D2").Value = 9.8857702241529104

$ws.Range("D9").Select()
